# Adding more home page testcases
# Populates the previously-empty "Schedule meeting from Home" (row 12) and
# "Join meeting from Home" (row 13) test cases with repro steps / expected
# results, and moves the active selection to D13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 12 - "[Zoom-Home] Schedule meeting from Home" (Sr no 9)
# ---------------------------------------------------------------------

$c12 = "Pre-requisite: user login`n1.In the right rail, click on Schedule present.`n2. Enter the meeting details and click Save.`n3. Go back to Home and observe the behavior."
$ws.Range("C12").Value = $c12

$c12Bold = $ws.Range("C12").Characters(1, 13)
$c12Bold.Font.Bold = $true
$c12Bold.Font.Name = "Calibri"
$c12Bold.Font.Size = 11

$c12Rest = $ws.Range("C12").Characters(14, ($c12.Length - 13))
$c12Rest.Font.Bold = $false
$c12Rest.Font.Name = "Calibri"
$c12Rest.Font.Size = 11

$d12 = "1. The top card should contain details like Schedule, Join and Host meetings and Personal meeting ID.`n1.1 Clicking on Schedule, should open PERSONAL -> Meetings page to schedule meeting in a new tab.`n2. Meeting should be created.`n3. Meeting card present in right rail, Home page should contain the scheduled meeting info."
$ws.Range("D12").Value = $d12

$ws.Rows.Item(12).RowHeight = 102

# ---------------------------------------------------------------------
# Row 13 - "[Zoom-Home] Join meeting from Home" (Sr no 10)
# ---------------------------------------------------------------------

$c13 = "Pre-requisite: user login`n1. In the right rail, click on 'Join' present.`n2. Enter the meeting ID/ link name or personal meeting id to start the meeting.`n3. Click on Join.`n4. A pop-up to open Zoom desktop app should show, click on Cancel.`n5. Click on Launch Meeting button present.`n6. Repeat step 4.`n7. Click on 'Join from your browser'..`n8. Obseve the behavior."
$ws.Range("C13").Value = $c13

$c13Bold = $ws.Range("C13").Characters(1, 13)
$c13Bold.Font.Bold = $true
$c13Bold.Font.Name = "Calibri"
$c13Bold.Font.Size = 11

$c13Rest = $ws.Range("C13").Characters(14, ($c13.Length - 13))
$c13Rest.Font.Bold = $false
$c13Rest.Font.Name = "Calibri"
$c13Rest.Font.Size = 11

$ws.Rows.Item(13).RowHeight = 145.75

# ---------------------------------------------------------------------
# Move the active selection to D13
# ---------------------------------------------------------------------

$ws.Range("D13").Select()
